{"js": "// Update the date line and the 25 division-practice answers in the table.\n// The document body starts with a centered date paragraph, followed by a\n// single 20-row x 5-column table where every 4th row (0, 4, 8, 12, 16)\n// holds the answer text and the rows in-between are blank spacer rows.\n\n// 1) Update the date heading paragraph.\nconst dateResults = context.document.body.search(\"2025-10-15 Wednesday\", { matchCase: true });\ndateResults.load(\"items\");\nawait context.sync();\nif (dateResults.items.length > 0) {\n  dateResults.items[0].insertText(\"2025-10-16 Thursday\", \"Replace\");\n}\n\n// 2) Update every answer cell in the table in one shot via Table.values,\n//    preserving the existing (blank) spacer rows untouched.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.values = [\n  [\"24\u00f74=6, 0\", \"97\u00f73=32, 1\", \"33\u00f78=4, 1\", \"56\u00f79=6, 2\", \"35\u00f76=5, 5\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"74\u00f73=24, 2\", \"22\u00f75=4, 2\", \"76\u00f78=9, 4\", \"22\u00f72=11, 0\", \"19\u00f75=3, 4\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"76\u00f74=19, 0\", \"40\u00f73=13, 1\", \"32\u00f79=3, 5\", \"50\u00f78=6, 2\", \"55\u00f73=18, 1\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"90\u00f79=10, 0\", \"22\u00f74=5, 2\", \"67\u00f76=11, 1\", \"67\u00f73=22, 1\", \"15\u00f79=1, 6\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"31\u00f75=6, 1\", \"27\u00f77=3, 6\", \"69\u00f77=9, 6\", \"19\u00f72=9, 1\", \"61\u00f77=8, 5\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n  [\"\", \"\", \"\", \"\", \"\"],\n];\n\nawait context.sync();\n", "ps1": "# Update the date line and the 25 division-practice answers in the table.\n# The document body starts with a centered date paragraph, followed by a\n# single 20-row x 5-column table where every 4th row (1, 5, 9, 13, 17 in\n# 1-based COM indexing) holds the answer text and the rows in-between are\n# blank spacer rows.\n\n$d = $word.ActiveDocument\n\n# 1) Update the date heading paragraph (first paragraph in the body).\n$d.Paragraphs.Item(1).Range.Text = \"2025-10-16 Thursday\"\n\n# 2) Update every answer cell in the table by (row, column) position so the\n#    blank spacer rows are left untouched and formatting on each cell is\n#    preserved.\n$t = $d.Tables.Item(1)\n\n$t.Cell(1,1).Range.Text = \"24\u00f74=6, 0\"\n$t.Cell(1,2).Range.Text = \"97\u00f73=32, 1\"\n$t.Cell(1,3).Range.Text = \"33\u00f78=4, 1\"\n$t.Cell(1,4).Range.Text = \"56\u00f79=6, 2\"\n$t.Cell(1,5).Range.Text = \"35\u00f76=5, 5\"\n\n$t.Cell(5,1).Range.Text = \"74\u00f73=24, 2\"\n$t.Cell(5,2).Range.Text = \"22\u00f75=4, 2\"\n$t.Cell(5,3).Range.Text = \"76\u00f78=9, 4\"\n$t.Cell(5,4).Range.Text = \"22\u00f72=11, 0\"\n$t.Cell(5,5).Range.Text = \"19\u00f75=3, 4\"\n\n$t.Cell(9,1).Range.Text = \"76\u00f74=19, 0\"\n$t.Cell(9,2).Range.Text = \"40\u00f73=13, 1\"\n$t.Cell(9,3).Range.Text = \"32\u00f79=3, 5\"\n$t.Cell(9,4).Range.Text = \"50\u00f78=6, 2\"\n$t.Cell(9,5).Range.Text = \"55\u00f73=18, 1\"\n\n$t.Cell(13,1).Range.Text = \"90\u00f79=10, 0\"\n$t.Cell(13,2).Range.Text = \"22\u00f74=5, 2\"\n$t.Cell(13,3).Range.Text = \"67\u00f76=11, 1\"\n$t.Cell(13,4).Range.Text = \"67\u00f73=22, 1\"\n$t.Cell(13,5).Range.Text = \"15\u00f79=1, 6\"\n\n$t.Cell(17,1).Range.Text = \"31\u00f75=6, 1\"\n$t.Cell(17,2).Range.Text = \"27\u00f77=3, 6\"\n$t.Cell(17,3).Range.Text = \"69\u00f77=9, 6\"\n$t.Cell(17,4).Range.Text = \"19\u00f72=9, 1\"\n$t.Cell(17,5).Range.Text = \"61\u00f77=8, 5\"\n"}
